# "Online Learning - Lihat Modul"
# Rename the lone sheet to "search", add four more sheets ("next page",
# "lihat modul", "materi belajar", "review"), fill the three data sheets
# with the module/search-result tables, and leave "materi belajar" as the
# active tab/sheet (matching the source workbook's final UI state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheets: rename existing Sheet1, then add the rest in an order that
#    reproduces the original sheetId allocation (1, 5, 2, 4, 3) once they
#    are rearranged into their final left-to-right order.
# ---------------------------------------------------------------------
$search = $wb.Worksheets.Item(1)
$search.Name = "search"

$lihatModul = $wb.Worksheets.Add($null, $search)
$lihatModul.Name = "lihat modul"

$review = $wb.Worksheets.Add($null, $lihatModul)
$review.Name = "review"

$materiBelajar = $wb.Worksheets.Add($null, $review)
$materiBelajar.Name = "materi belajar"

$nextPage = $wb.Worksheets.Add($null, $materiBelajar)
$nextPage.Name = "next page"

# Re-fetch by name before each Move (stable across reindexing).
$wb.Worksheets.Item("next page").Move($null, $wb.Worksheets.Item("search"))
$wb.Worksheets.Item("materi belajar").Move($wb.Worksheets.Item("review"), $null)

# Final left-to-right order: search, next page, lihat modul, materi belajar, review

# ---------------------------------------------------------------------
# 2. "lihat modul" data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("lihat modul")
$ws.Range("A1").Value = "module_name"
$ws.Range("B1").Value = "found"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "expected"

$ws.Range("A2").Value = "zzzzzzzzzz"
$ws.Range("B2").Value = "no"
$ws.Range("D2").Value = "fail"

$ws.Range("A3").Value = "WELCOMING KIT"
$ws.Range("B3").Value = "yes"
$ws.Range("D3").Value = "pass"

$ws.Range("A4").Value = "NEOP AAV"
$ws.Range("B4").Value = "yes"
$ws.Range("C4").Value = "new"
$ws.Range("D4").Value = "pass"

$ws.Columns.Item(1).AutoFit()
$ws.Range("A3").Select()

# ---------------------------------------------------------------------
# 3. "materi belajar" data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("materi belajar")
$ws.Range("A1").Value = "module_name"
$ws.Range("B1").Value = "found"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "expected"

$ws.Range("A2").Value = "WELCOMING KIT"
$ws.Range("B2").Value = "yes"
$ws.Range("D2").Value = "pass"

$ws.Columns.Item(1).AutoFit()
$ws.Range("L16").Select()

# ---------------------------------------------------------------------
# 4. "review" data
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("review")
$ws.Range("A1").Value = "module_name"
$ws.Range("B1").Value = "found"
$ws.Range("C1").Value = "status"
$ws.Range("D1").Value = "expected"

$ws.Range("A2").Value = "yyyyyyyyyyy"
$ws.Range("B2").Value = "no"
$ws.Range("D2").Value = "fail"

$ws.Range("A3").Value = "NEOP AAV"
$ws.Range("B3").Value = "yes"
$ws.Range("C3").Value = "review empty"
$ws.Range("D3").Value = "fail"

$ws.Range("A4").Value = "WELCOMING KIT"
$ws.Range("B4").Value = "yes"

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Range("A3").Select()

# ---------------------------------------------------------------------
# 5. "search" sheet: clear the tab selection / set new cell selection
#    (data itself is unchanged).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("search")
$ws.Range("A3").Select()

# ---------------------------------------------------------------------
# 6. Active sheet/tab ends on "materi belajar" (activeTab=3, 0-based).
# ---------------------------------------------------------------------
$materiBelajarFinal = $wb.Worksheets.Item("materi belajar")
$materiBelajarFinal.Select()
$materiBelajarFinal.Range("L16").Select()
